$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (2021-04-27 .. 2021-05-02), appended after the existing
# last row (238, date serial 44312).
$data = @(
    @(239, 44313, 15, 68,  93.96288466055908),
    @(240, 44314,  2, 68,  93.96288466055908),
    @(241, 44315,  8, 71,  98.10830604264255),
    @(242, 44316, 32, 92, 127.126255717227),
    @(243, 44317, 16, 94, 129.8898699719493),
    @(244, 44318, 14, 94, 129.8898699719493)
)

foreach ($entry in $data) {
    $row = $entry[0]

    # Column A carries the same date style ("s=2") as every prior row in
    # the column; copy it from the row directly above before overwriting
    # the value so formatting matches the rest of the series.
    $ws.Range("A" + ($row - 1)).Copy($ws.Range("A" + $row))

    $ws.Cells.Item($row, 1).Value = $entry[1]
    $ws.Cells.Item($row, 2).Value = $entry[2]
    $ws.Cells.Item($row, 3).Value = $entry[3]
    $ws.Cells.Item($row, 4).Value = $entry[4]
}
